$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unmerge the old title cell and clear all existing content ---
$ws.Range("C3:F3").UnMerge()
$ws.Cells.Clear()

# --- Column widths ---
$ws.Columns("A").ColumnWidth = 0.5
$ws.Range("B1:G1").EntireColumn.ColumnWidth = 30.88671875
$ws.Range("H1:L1").EntireColumn.ColumnWidth = 14.109375

# --- Row 3: title, merged B3:G3 ---
$ws.Range("B3").Value = "Teste de Mesa - Fração"
$ws.Range("B3:G3").HorizontalAlignment = -4108
$ws.Range("B3:G3").VerticalAlignment = -4108
$ws.Range("H3:L3").VerticalAlignment = -4108
$ws.Range("B3:G3").Merge()

# --- Row 4: headers ---
$ws.Range("B4").Value = "Numerador Primeira Fracao"
$ws.Range("C4").Value = "Denominador Primeira Fracao"
$ws.Range("D4").Value = "Numerador Segunda Fracao"
$ws.Range("E4").Value = "Denominador Segunda Fracao"
$ws.Range("F4").Value = "Equação"
$ws.Range("G4").Value = "Resultado"
$ws.Range("B4:G4").HorizontalAlignment = -4108
$ws.Range("B4:G4").VerticalAlignment = -4108

# --- Row 5: Divisão => 5/4 ÷ 1/3 = 15/4 ---
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = "Divisão"
$r = $ws.Range("G5")
$r.Value = "15/4*"
$r.Characters(5,1).Font.Color = 16777215

# --- Row 6: Soma => 2/5 + 3/7 = 29/35 ---
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = "Soma"
$r = $ws.Range("G6")
$r.Value = "29/35*"
$r.Characters(6,1).Font.Color = 16777215

# --- Row 7: Subtração => 10/3 - 4/3 = 2 ---
$ws.Range("B7").Value = 10
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 4
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = "Subtração"
$ws.Range("G7").Value = 2

# --- Row 8: Multiplicação => 5/3 * 7/4 = 35/12 ---
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = 3
$ws.Range("D8").Value = 7
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = "Multiplicação"
$r = $ws.Range("G8")
$r.Value = "35/12*"
$r.Characters(6,1).Font.Color = 16777215

# --- Alignment for the data block (rows 5-8, B:G) to match header style ---
$ws.Range("B5:G8").HorizontalAlignment = -4108
$ws.Range("B5:G8").VerticalAlignment = -4108

# --- View: selection + scroll position ---
$ws.Range("G13").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 1
